$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("模块名")

# The 6th "round" result column header is replaced with a generic
# "测试结果" (Test Result) header, styled in SimSun like the sheet's
# other Chinese-font headers (e.g. B2 "所属模块").
$ws.Range("K1").Value = "测试结果"
$ws.Range("K1").Font.Name = "宋体"

# Selection moves from the old O3 to K2, matching the saved view state.
$ws.Range("K2").Select()
